$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.010.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.17%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.299.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.51%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'301.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.21%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'98.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.27%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.68%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.505"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.13%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'34.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.39%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +4.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'17.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +15.75%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.656.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.56%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.293.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.31%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.809"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.76%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'42.922.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.99%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +8.75%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.91%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.28%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.52%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'236.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.46%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +14.79%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.12%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.64%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'24.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.96%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'167.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -9.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'33.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.78%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'9.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'WEMIXToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.19%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'RenderToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'16.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.65%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E39").Value = "'  +3.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.88%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -5.28%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.990.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.60%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.59%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'17.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.37%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'56.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +9.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.524.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.10%  "
$ws.Range("E51").Style = "Normal"
